$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.981.42'
$ws.Range('E2').Value = '  +0.48%  '

$ws.Range('D3').Value = '1.641.01'
$ws.Range('E3').Value = '  +0.59%  '

$ws.Range('E4').Value = '  +0.30%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.77'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E6').Value = '  +0.23%  '

$ws.Range('E7').Value = '  +0.29%  '

$ws.Range('E8').Value = '  +0.14%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0638'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.99%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.55'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.13%  '

$ws.Range('E11').Value = '  +0.46%  '

$ws.Range('D12').Value = '1.869.24'
$ws.Range('E12').Value = '  +0.58%  '

$ws.Range('D13').Value = '1.669.42'
$ws.Range('E13').Value = '  +0.18%  '

$ws.Range('E14').Value = '  +0.20%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.544'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.13%  '

$ws.Range('E16').Value = '  +1.04%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.32'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.11%  '

$ws.Range('D18').Value = '26.083.52'
$ws.Range('E18').Value = '  +0.81%  '

$ws.Range('E19').Value = '  +0.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.97%  '

$ws.Range('E21').Value = '  -0.58%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.92'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.02%  '

$ws.Range('E23').Value = '  -0.68%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.132'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.31%  '

$ws.Range('E25').Value = '  -2.36%  '

$ws.Range('E26').Value = '  +0.54%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '142.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.21%  '

$ws.Range('E28').Value = '  +0.61%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.70%  '

$ws.Range('E30').Value = '  +0.85%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0496'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.34%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.29'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.07%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.25'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.00%  '

$ws.Range('E34').Value = '  -1.36%  '

$ws.Range('E35').Value = '  +1.31%  '

$ws.Range('E36').Value = '  +0.44%  '

$ws.Range('D37').Value = '1.124.45'
$ws.Range('E37').Value = '  -1.19%  '

$ws.Range('E38').Value = '  -1.25%  '

$ws.Range('E39').Value = '  -0.57%  '

$ws.Range('E40').Value = '  +0.27%  '

$ws.Range('E41').Value = '  +0.54%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.07'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.00%  '

$ws.Range('E43').Value = '  +0.13%  '

$ws.Range('D44').Value = '1.778.27'
$ws.Range('E44').Value = '  +0.60%  '

$ws.Range('D45').Value = '0.0₆0117'
$ws.Range('E45').Value = '  +3.67%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.32'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.14%  '

$ws.Range('E47').Value = '  -0.62%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.48'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.96%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.75'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.04%  '

$ws.Range('E50').Value = '  -0.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0953'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.65%  '
